# "removed old backup, replaced switch with jumper"
# The BOM's LED (D2, row 5) LCSC part number changes from the old/backup
# part C401117 to the replacement part C409784.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "C409784"
